# Updated test data for 5,24,40V,BatteryStandby and AC Calculations test cases
#
# "Add Panels" sheet (sheet1): the 40V load row (row 8) had its CPU Type
# (C8) cleared out and its two 40V-rail readings (J8, N8) bumped from
# "0.00" to "0.000". The four extra sample rows (9-12: MZX 125, Pro32xD,
# MX2-100, P885D) are no longer needed here - three of them (MZX 125,
# Pro32xD, P885D) move over to the "Test data" sheet, and the MX2-100 row
# is dropped entirely.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Add Panels"
$ws2 = $wb.Worksheets.Item(2)   # "Test data"

# --- Move the still-needed rows over to the "Test data" sheet first,
#     before anything on "Add Panels" gets deleted/renumbered. ---
$ws1.Range("A9:F9").Copy($ws2.Range("A1:F1"))     # MZX 125
$ws1.Range("A10:F10").Copy($ws2.Range("A2:F2"))   # Pro32xD
# (row 11, MX2-100, is intentionally skipped - it is not carried over)
$ws1.Range("A12:F12").Copy($ws2.Range("A3:F3"))   # P885D

# The copied "CPU Type" cells need to show as blank text (quote-prefixed
# empty-text style, same as used elsewhere in this column) - re-stamp all
# three, since a literal empty value doesn't round-trip as a typed cell.
$ws2.Range("C1").Formula = "'"
$ws2.Range("C2").Formula = "'"
$ws2.Range("C3").Formula = "'"

$ws2.Range("A1:XFD3").Select()

# --- Clean up "Add Panels": clear the CPU Type on the 40V load row and
#     update its two readings, then drop the now-relocated sample rows. ---
$ws1.Range("C8").Formula = "'"
$ws1.Range("J8").Formula = "'0.000"
$ws1.Range("N8").Formula = "'0.000"

$ws1.Range("A9:N12").EntireRow.Delete()

$ws1.Activate()
$ws1.Range("N8").Select()
